$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18-20 down to 19-21)
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 (Colombia Primera A: Atl. Nacional vs Santa Fe)
$ws.Range("A18").Value = "E1chGh3C"
$ws.Range("B18").Value = "20/11/2024"
$ws.Range("C18").Value = "22:30"
$ws.Range("D18").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E18").Value = "Atl. Nacional"
$ws.Range("F18").Value = "Santa Fe"
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 4.1
$ws.Range("J18").Value = 2.75
$ws.Range("K18").Value = 1.95
$ws.Range("L18").Value = 4.75
$ws.Range("M18").Value = 1.08
$ws.Range("N18").Value = 8
$ws.Range("O18").Value = 1.44
$ws.Range("P18").Value = 2.63
$ws.Range("Q18").Value = 2.35
$ws.Range("R18").Value = 1.57
$ws.Range("S18").Value = 1.53
$ws.Range("T18").Value = 2.38
$ws.Range("U18").Value = 2.1
$ws.Range("V18").Value = 1.67
$ws.Range("W18").Value = 5.5
$ws.Range("X18").Value = 8
$ws.Range("Y18").Value = 9.5
$ws.Range("Z18").Value = 17
$ws.Range("AA18").Value = 19
$ws.Range("AB18").Value = 41
$ws.Range("AC18").Value = 7
$ws.Range("AD18").Value = 6.5
$ws.Range("AE18").Value = 19
$ws.Range("AF18").Value = 67
$ws.Range("AG18").Value = 9
$ws.Range("AH18").Value = 19
$ws.Range("AI18").Value = 15
$ws.Range("AJ18").Value = 41
$ws.Range("AK18").Value = 41
$ws.Range("AL18").Value = 51
$ws.Range("AM18").Value = 1000
$ws.Range("AN18").Value = 3.75
$ws.Range("AO18").Value = 11
$ws.Range("AP18").Value = 26
$ws.Range("AQ18").Value = 41
$ws.Range("AR18").Value = 67
$ws.Range("AS18").Value = 251
$ws.Range("AT18").Value = 2.38
$ws.Range("AU18").Value = 9.5
$ws.Range("AV18").Value = 67
$ws.Range("AW18").Value = 6
$ws.Range("AX18").Value = 26
$ws.Range("AY18").Value = 41
$ws.Range("AZ18").Value = 81
$ws.Range("BA18").Value = 126
$ws.Range("BB18").Value = 351
$ws.Range("BC18").Value = 126
$ws.Range("BD18").Value = 126

# Update odds values in other rows that changed independently of the row insertion
$ws.Range("G6").Value = 2
$ws.Range("G9").Value = 1.48
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 2.05
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 7
$ws.Range("Q9").Value = 2.06
$ws.Range("R9").Value = 1.84
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("X9").Value = 6.5
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 9.5
$ws.Range("AD9").Value = 8
$ws.Range("AE9").Value = 21
$ws.Range("AN9").Value = 3.25
$ws.Range("AO9").Value = 7.5
$ws.Range("AQ9").Value = 23
$ws.Range("AU9").Value = 9.5
$ws.Range("AW9").Value = 8
$ws.Range("AZ9").Value = 151
$ws.Range("BB9").Value = 401
$ws.Range("Q10").Value = 2.06
$ws.Range("R10").Value = 1.84
$ws.Range("G11").Value = 2.63
$ws.Range("J11").Value = 3.4
$ws.Range("K11").Value = 1.91
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 1.5
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("U11").Value = 2.05
$ws.Range("V11").Value = 1.7
$ws.Range("W11").Value = 6.5
$ws.Range("Y11").Value = 11
$ws.Range("AA11").Value = 26
$ws.Range("AC11").Value = 6.5
$ws.Range("AM11").Value = 501
$ws.Range("AS11").Value = 301
$ws.Range("AT11").Value = 2.25
$ws.Range("AY11").Value = 34
$ws.Range("G13").Value = 3.2
$ws.Range("I13").Value = 2.2
$ws.Range("L13").Value = 2.88
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4.33
$ws.Range("AH13").Value = 12
$ws.Range("AN13").Value = 5
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 2.75
$ws.Range("Q14").Value = 2.3
$ws.Range("R14").Value = 1.6
$ws.Range("Q15").Value = 1.98
$ws.Range("R15").Value = 1.92
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5

# Update the Wales match row (now row 21 after the insert) with its revised odds
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 8.25
$ws.Range("I21").Value = 1.06
$ws.Range("J21").Value = 21
$ws.Range("K21").Value = 3.45
$ws.Range("L21").Value = 1.29
$ws.Range("O21").Value = 1.07
$ws.Range("P21").Value = 6.8
$ws.Range("Q21").Value = 1.24
$ws.Range("R21").Value = 3.6
$ws.Range("S21").Value = 1.17
$ws.Range("T21").Value = 4.4
$ws.Range("U21").Value = 2.4
$ws.Range("V21").Value = 1.5
$ws.Range("W21").Value = 175
$ws.Range("X21").Value = 900
$ws.Range("Y21").Value = 175
$ws.Range("AB21").Value = 500
$ws.Range("AC21").Value = 25
$ws.Range("AD21").Value = 26
$ws.Range("AE21").Value = 55
$ws.Range("AF21").Value = 250
$ws.Range("AG21").Value = 12.5
$ws.Range("AH21").Value = 7.4
$ws.Range("AI21").Value = 14.5
$ws.Range("AK21").Value = 12
$ws.Range("AL21").Value = 45
$ws.Range("AN21").Value = 30
$ws.Range("AO21").Value = 300
$ws.Range("AP21").Value = 150
$ws.Range("AT21").Value = 4.4
$ws.Range("AU21").Value = 12.5
$ws.Range("AV21").Value = 110
$ws.Range("AX21").Value = 3.95
$ws.Range("AY21").Value = 15
$ws.Range("AZ21").Value = 7
$ws.Range("BB21").Value = 200
